$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ idx=4; D=44459; L="Especial"; M=200; N=2600; O=2700; P=2650; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2650; T=1 },
  @{ idx=5; D=44459; L="Primera"; M=300; N=2200; O=2300; P=2250; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2250; T=1 },
  @{ idx=6; D=44459; L="Segunda"; M=240; N=1900; O=2000; P=1950; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1950; T=1 },
  @{ idx=7; D=44169; L="Especial"; M=240; N=14500; O=15000; P=14750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1844; T=8 },
  @{ idx=8; D=44169; L="Primera"; M=240; N=12500; O=13000; P=12750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1594; T=8 },
  @{ idx=9; D=44169; L="Segunda"; M=200; N=10000; O=10500; P=10250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1281; T=8 },
  @{ idx=10; D=44449; L="Especial"; M=240; N=2900; O=3000; P=2950; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2950; T=1 },
  @{ idx=11; D=44449; L="Extra (doble especial)"; M=160; N=3100; O=3200; P=3150; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=3150; T=1 },
  @{ idx=12; D=44449; L="Primera"; M=300; N=2700; O=2800; P=2750; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2750; T=1 },
  @{ idx=13; D=44161; L="Especial"; M=240; N=13500; O=14000; P=13750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1719; T=8 },
  @{ idx=14; D=44161; L="Primera"; M=300; N=11500; O=12000; P=11750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1469; T=8 },
  @{ idx=15; D=44161; L="Segunda"; M=200; N=9000; O=9500; P=9250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1156; T=8 },
  @{ idx=16; D=44165; L="Especial"; M=300; N=14000; O=14500; P=14250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1781; T=8 },
  @{ idx=17; D=44165; L="Primera"; M=240; N=12000; O=12500; P=12250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1531; T=8 },
  @{ idx=18; D=44165; L="Segunda"; M=200; N=9500; O=10000; P=9750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1219; T=8 },
  @{ idx=19; D=44172; L="Especial"; M=200; N=14000; O=14500; P=14250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1781; T=8 },
  @{ idx=20; D=44172; L="Primera"; M=200; N=12000; O=12500; P=12250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1531; T=8 },
  @{ idx=21; D=44172; L="Segunda"; M=200; N=9500; O=10000; P=9750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1219; T=8 },
  @{ idx=22; D=44166; L="Especial"; M=300; N=14000; O=14500; P=14250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1781; T=8 },
  @{ idx=23; D=44166; L="Primera"; M=200; N=12000; O=12500; P=12250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1531; T=8 },
  @{ idx=24; D=44168; L="Especial"; M=240; N=14000; O=14500; P=14250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1781; T=8 },
  @{ idx=25; D=44168; L="Primera"; M=200; N=12000; O=12500; P=12250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1531; T=8 },
  @{ idx=26; D=44168; L="Segunda"; M=200; N=9500; O=10000; P=9750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1219; T=8 },
  @{ idx=27; D=44162; L="Especial"; M=340; N=14000; O=14500; P=14250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1781; T=8 },
  @{ idx=28; D=44162; L="Primera"; M=300; N=12000; O=12500; P=12250; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1531; T=8 },
  @{ idx=29; D=44162; L="Segunda"; M=200; N=9500; O=10000; P=9750; Q="`$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1219; T=8 },
  @{ idx=30; D=44410; L="Primera"; M=240; N=2400; O=2500; P=2450; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia del Elquí"; S=2450; T=1 },
  @{ idx=31; D=44410; L="Segunda"; M=240; N=2000; O=2100; P=2050; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia del Elquí"; S=2050; T=1 },
  @{ idx=32; D=44410; L="Tercera"; M=200; N=1600; O=1700; P=1650; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia del Elquí"; S=1650; T=1 },
  @{ idx=33; D=44411; L="Primera"; M=600; N=2400; O=2500; P=2450; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia del Elquí"; S=2450; T=1 },
  @{ idx=34; D=44411; L="Segunda"; M=400; N=2000; O=2100; P=2050; Q="`$/kilo (en caja de 15 kilos)"; R="Provincia del Elquí"; S=2050; T=1 }
)

$dateFormat = $ws.Range("D2").NumberFormat()

foreach ($row in $rows) {
    $i = $row.idx
    if ($i -gt 31) {
        $ws.Range("A" + $i).Value = 8
        $ws.Range("B" + $i).Value = "Terminal La Palmera de La Serena"
        $ws.Range("C" + $i).Value = "Coquimbo"
        $ws.Range("E" + $i).Value = 4
        $ws.Range("F" + $i).Value = "Fruta"
        $ws.Range("G" + $i).Value = 100107
        $ws.Range("H" + $i).Value = "Otros"
        $ws.Range("I" + $i).Value = 100107002
        $ws.Range("J" + $i).Value = "Chirimoya"
        $ws.Range("K" + $i).Value = "Cultivar IV Región"
    }
    $ws.Range("D" + $i).Value = $row.D
    $ws.Range("D" + $i).NumberFormat = $dateFormat
    $ws.Range("L" + $i).Value = $row.L
    $ws.Range("M" + $i).Value = $row.M
    $ws.Range("N" + $i).Value = $row.N
    $ws.Range("O" + $i).Value = $row.O
    $ws.Range("P" + $i).Value = $row.P
    $ws.Range("Q" + $i).Value = $row.Q
    $ws.Range("R" + $i).Value = $row.R
    $ws.Range("S" + $i).Value = $row.S
    $ws.Range("T" + $i).Value = $row.T
}